$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lernjournal Aufgabe 3")

# New journal entry text (used twice, shared string)
$entryText = "Diagramme überarbeitet und Glossar verfeinert." + [char]10 + "Abgabefertig für's erste."

# --- Section 1 (rows 11-24): fill empty row 19 ---
$ws.Range("A19").Value = $entryText
$ws.Range("B19").Value = 45
$ws.Range("C19").Value = 41545
$ws.Rows.Item(19).RowHeight = 24

# --- Section 2 (rows 31-45): fill empty row 41, matching formatting of row 19 ---
$row41 = $ws.Range("A41:C41")
$row41.Font.Italic = $false
$row41.Borders.LineStyle = 1
$row41.Borders.Weight = 1

$ws.Range("A41").Value = $entryText
$ws.Range("B41").Value = 45
$ws.Range("C41").Value = 41545
$ws.Range("C41").NumberFormat = $ws.Range("C19").NumberFormat
$ws.Rows.Item(41).RowHeight = 24

$excel.CutCopyMode = 0

# --- Update sheet view / selection to match author's saved state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$ws.Range("C58").Select()

$wb.Save()
